$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1, mirrors the style of the other header cells (row 1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# New data values for K2:K7
$ws.Range("K2").Value = "PROCEDURE"
$ws.Range("K3").Value = "BEHAVIORAL"
$ws.Range("K4").Value = "OTHER"
$ws.Range("K5").Value = "OTHER"
$ws.Range("K6").Value = "OTHER"
$ws.Range("K7").Value = "BEHAVIORAL"
